$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '57.828.92'
$ws.Cells.Item(2, 5).Value = '  -5.46%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.913.51'
$ws.Cells.Item(3, 5).Value = '  -3.26%  '
$ws.Cells.Item(4, 5).Value = '  +0.09%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '551.22'
$ws.Cells.Item(5, 5).Value = '  -2.99%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '123.70'
$ws.Cells.Item(6, 5).Value = '  -4.24%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '2.912.67'
$ws.Cells.Item(8, 5).Value = '  -3.25%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.497'
$ws.Cells.Item(9, 5).Value = '  +0.05%  '
$ws.Cells.Item(10, 5).Value = '  -6.33%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '4.81'
$ws.Cells.Item(11, 5).Value = '  -7.38%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.442'
$ws.Cells.Item(12, 5).Value = '  +2.86%  '
$ws.Cells.Item(13, 5).Value = '  -4.50%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '32.50'
$ws.Cells.Item(14, 5).Value = '  -1.17%  '
$ws.Cells.Item(15, 5).Value = '  +1.30%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '3.393.46'
$ws.Cells.Item(16, 5).Value = '  -3.07%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.914.13'
$ws.Cells.Item(17, 5).Value = '  -3.07%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '6.62'
$ws.Cells.Item(18, 5).Value = '  +6.48%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '57.838.02'
$ws.Cells.Item(19, 5).Value = '  -5.51%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '411.11'
$ws.Cells.Item(20, 5).Value = '  -6.52%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '12.93'
$ws.Cells.Item(21, 5).Value = '  -2.01%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.676'
$ws.Cells.Item(22, 5).Value = '  +2.01%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '6.88'
$ws.Cells.Item(23, 5).Value = '  -3.72%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '12.93'
$ws.Cells.Item(24, 5).Value = '  +3.05%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '77.37'
$ws.Cells.Item(25, 5).Value = '  -1.99%  '
$ws.Cells.Item(26, 5).Value = '  +0.13%  '
$ws.Cells.Item(27, 5).Value = '  +0.04%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.47'
$ws.Cells.Item(28, 5).Value = '  -0.86%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '7.30'
$ws.Cells.Item(29, 5).Value = '  +0.98%  '
$ws.Cells.Item(30, 5).Value = '  +3.73%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '6.10'
$ws.Cells.Item(31, 5).Value = '  -1.85%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '24.76'
$ws.Cells.Item(32, 5).Value = '  -2.98%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.0986'
$ws.Cells.Item(33, 5).Value = '  +4.65%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.913'
$ws.Cells.Item(34, 5).Value = '  -4.30%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '5.42'
$ws.Cells.Item(35, 5).Value = '  -2.74%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.02'
$ws.Cells.Item(36, 5).Value = '  -11.00%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '48.22'
$ws.Cells.Item(37, 5).Value = '  -3.72%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '8.45'
$ws.Cells.Item(38, 5).Value = '  +9.15%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.0₃0629'
$ws.Cells.Item(39, 5).Value = '  -7.89%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.0347'
$ws.Cells.Item(40, 5).Value = '  -4.60%  '
$ws.Cells.Item(41, 5).Value = '  -1.26%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '2.633.61'
$ws.Cells.Item(42, 5).Value = '  -0.60%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '364.34'
$ws.Cells.Item(43, 5).Value = '  -2.35%  '
$ws.Cells.Item(44, 5).Value = '  -0.78%  '
$ws.Cells.Item(45, 5).Value = '  +0.00%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '119.70'
$ws.Cells.Item(46, 5).Value = '  +0.16%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.230'
$ws.Cells.Item(47, 5).Value = '  -2.41%  '
$ws.Cells.Item(48, 2).Value = 'Stellar'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.107'
$ws.Cells.Item(48, 5).Value = '  +1.02%  '
$ws.Cells.Item(49, 2).Value = 'Fetch.AI'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.96'
$ws.Cells.Item(49, 5).Value = '  -0.31%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '22.85'
$ws.Cells.Item(50, 5).Value = '  -3.10%  '
$ws.Cells.Item(51, 5).Value = '  -2.50%  '
